$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.62
$ws.Range("H2").Value = 3.9
$ws.Range("J2").Value = 2.25
$ws.Range("K2").Value = 2.1
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 8.5
$ws.Range("U2").Value = 2.2
$ws.Range("V2").Value = 1.62
$ws.Range("W2").Value = 5.5
$ws.Range("X2").Value = 6.5
$ws.Range("Y2").Value = 9
$ws.Range("AC2").Value = 8.5
$ws.Range("AF2").Value = 81
$ws.Range("AH2").Value = 11
$ws.Range("AK2").Value = 51
$ws.Range("AP2").Value = 23
$ws.Range("AQ2").Value = 29
$ws.Range("AS2").Value = 201
$ws.Range("AU2").Value = 9.5

# Row 4 updates
$ws.Range("G4").Value = 3.1
$ws.Range("K4").Value = 2.2
$ws.Range("L4").Value = 2.62
$ws.Range("Q4").Value = 1.65
$ws.Range("R4").Value = 2
$ws.Range("AB4").Value = 28
$ws.Range("AC4").Value = 12.5
$ws.Range("AI4").Value = 11.5
$ws.Range("AJ4").Value = 8.5
$ws.Range("AN4").Value = 5.2
$ws.Range("AT4").Value = 3
$ws.Range("AU4").Value = 6.5
$ws.Range("AX4").Value = 10.5
$ws.Range("AY4").Value = 16.5
$ws.Range("AZ4").Value = 37
$ws.Range("BA4").Value = 60

$wb.Save()
